$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "'28.017.12"
$ws.Range("E2").Value = "  -0.24%  "

$ws.Range("D3").Formula = "'1.870.49"
$ws.Range("E3").Value = "  -0.21%  "

$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").Formula = "'311.83"
$ws.Range("E5").Value = "  -0.45%  "

$ws.Range("E6").Value = "  +0.10%  "

$ws.Range("D7").Formula = "'0.5158"
$ws.Range("E7").Value = "  +2.17%  "

$ws.Range("D8").Formula = "'0.3856"
$ws.Range("E8").Value = "  +0.43%  "

$ws.Range("E9").Value = "  -3.23%  "

$ws.Range("D10").Formula = "'1.111"
$ws.Range("E10").Value = "  -0.27%  "

$ws.Range("D11").Formula = "'41.52"
$ws.Range("E11").Value = "  +0.55%  "

$ws.Range("D12").Formula = "'6.208"
$ws.Range("E12").Value = "  -1.34%  "

$ws.Range("D13").Formula = "'20.56"
$ws.Range("E13").Value = "  -0.42%  "

$ws.Range("D14").Formula = "'1.858.01"
$ws.Range("E14").Value = "  -1.19%  "

$ws.Range("D15").Formula = "'7.317"
$ws.Range("E15").Value = "  +1.52%  "

$ws.Range("D16").Formula = "'1.003"
$ws.Range("E16").Value = "  -0.01%  "

$ws.Range("D17").Formula = "'0.00001098"
$ws.Range("E17").Value = "  +0.01%  "

$ws.Range("D18").Formula = "'90.85"
$ws.Range("E18").Value = "  -0.17%  "

$ws.Range("D19").Formula = "'0.06639"
$ws.Range("E19").Value = "  +0.17%  "

$ws.Range("D20").Formula = "'17.72"
$ws.Range("E20").Value = "  -1.98%  "

$ws.Range("E21").Value = "  +0.11%  "

$ws.Range("D22").Formula = "'6.039"
$ws.Range("E22").Value = "  -0.89%  "

$ws.Range("D23").Formula = "'28.059.12"
$ws.Range("E23").Value = "  -0.24%  "

$ws.Range("D24").Formula = "'11.12"
$ws.Range("E24").Value = "  -2.22%  "

$ws.Range("E25").Value = "  -0.60%  "

$ws.Range("D26").Formula = "'2.089.73"
$ws.Range("E26").Value = "  -0.24%  "

$ws.Range("D27").Formula = "'2.506"
$ws.Range("E27").Value = "  -3.08%  "

$ws.Range("D28").Formula = "'157.48"
$ws.Range("E28").Value = "  +0.38%  "

$ws.Range("D29").Formula = "'20.66"
$ws.Range("E29").Value = "  -0.25%  "

$ws.Range("D30").Formula = "'125.16"
$ws.Range("E30").Value = "  -0.90%  "

$ws.Range("D31").Formula = "'0.1066"
$ws.Range("E31").Value = "  +0.81%  "

$ws.Range("D32").Formula = "'1.033"
$ws.Range("E32").Value = "  -2.77%  "

$ws.Range("D33").Formula = "'5.797"
$ws.Range("E33").Value = "  +3.16%  "

$ws.Range("D34").Formula = "'3.602"
$ws.Range("E34").Value = "  +0.49%  "

$ws.Range("D35").Formula = "'9.482"
$ws.Range("E35").Value = "  -1.46%  "

$ws.Range("E36").Value = "  -1.04%  "

$ws.Range("D37").Formula = "'0.06510"
$ws.Range("E37").Value = "  -1.10%  "

$ws.Range("D38").Formula = "'0.2203"
$ws.Range("E38").Value = "  +1.03%  "

$ws.Range("D39").Formula = "'0.6588"
$ws.Range("E39").Value = "  +3.30%  "

$ws.Range("D40").Formula = "'1.202"
$ws.Range("E40").Value = "  -0.94%  "

$ws.Range("D41").Formula = "'5.010"
$ws.Range("E41").Value = "  +2.31%  "

$ws.Range("D42").Formula = "'1.209"
$ws.Range("E42").Value = "  -2.65%  "

$ws.Range("D43").Formula = "'11.20"
$ws.Range("E43").Value = "  -2.00%  "

$ws.Range("D44").Formula = "'0.6144"
$ws.Range("E44").Value = "  +2.31%  "

$ws.Range("D45").Formula = "'13.09"
$ws.Range("E45").Value = "  -0.65%  "

$ws.Range("E46").Value = "  +0.02%  "

$ws.Range("D47").Formula = "'3.670"

$ws.Range("D48").Formula = "'2.026"
$ws.Range("E48").Value = "  +1.86%  "

$ws.Range("D49").Formula = "'1.220"
$ws.Range("E49").Value = "  -0.26%  "

$ws.Range("D50").Formula = "'121.09"
$ws.Range("E50").Value = "  -0.28%  "

$ws.Range("E51").Value = "  -2.74%  "
